$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.484.19"
$ws.Range("E2").Value = "  +4.92%  "

$ws.Range("D3").Value = "2.721.41"
$ws.Range("E3").Value = "  +3.26%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'578.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").Value = "'153.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.87%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("E8").Value = "  +1.57%  "

$ws.Range("D9").Value = "2.746.68"
$ws.Range("E9").Value = "  +4.32%  "

$ws.Range("E10").Value = "  +3.11%  "

$ws.Range("E11").Value = "  +6.13%  "

$ws.Range("D12").Value = "'0.389"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.76%  "

$ws.Range("D14").Value = "3.208.37"
$ws.Range("E14").Value = "  +3.95%  "

$ws.Range("D15").Value = "'26.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("D16").Value = "63.377.38"
$ws.Range("E16").Value = "  +4.77%  "

$ws.Range("E17").Value = "  +6.84%  "

$ws.Range("D18").Value = "2.747.20"
$ws.Range("E18").Value = "  +4.58%  "

$ws.Range("D19").Value = "'11.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.46%  "

$ws.Range("D20").Value = "'4.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.95%  "

$ws.Range("D21").Value = "'359.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.45%  "

$ws.Range("D22").Value = "'6.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.30%  "

$ws.Range("D23").Value = "'0.539"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").Value = "'0.994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("D25").Value = "'65.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.51%  "

$ws.Range("E26").Value = "  +4.49%  "

$ws.Range("D27").Value = "'8.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.31%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("E29").Value = "  +12.09%  "

$ws.Range("E30").Value = "  -2.57%  "

$ws.Range("D31").Value = "'7.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.31%  "

$ws.Range("D32").Value = "'172.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.67%  "

$ws.Range("E33").Value = "  +12.65%  "

$ws.Range("E34").Value = "  -0.11%  "

$ws.Range("D35").Value = "'20.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.44%  "

$ws.Range("E36").Value = "  +7.39%  "

$ws.Range("E37").Value = "  +9.09%  "

$ws.Range("D38").Value = "'1.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.98%  "

$ws.Range("D39").Value = "'0.991"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +15.10%  "

$ws.Range("D40").Value = "'344.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.00%  "

$ws.Range("D41").Value = "'4.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.67%  "

$ws.Range("D42").Value = "'39.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("D43").Value = "'5.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.83%  "

$ws.Range("D44").Value = "'21.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.28%  "

$ws.Range("D45").Value = "'21.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.69%  "

$ws.Range("D46").Value = "'0.646"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.08%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'138.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.55%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0585"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.92%  "

$ws.Range("D49").Value = "'0.0255"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.42%  "

$ws.Range("E50").Value = "  +1.44%  "

$ws.Range("D51").Value = "'0.996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.41%  "
